$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace the single example download entry with the new one ---
$ws.Range("A2").Value = "S1-12proteomics.xlsx"
$ws.Range("B2").Value = "Download .xlsx"
$ws.Range("C2").Value = "Download the data collected in this study as an Excel spreadsheet"
$ws.Range("D2").Value = "S1-12"

# Row 2 no longer carries the wrap/valign style used by the data rows below
$ws.Range("A2:D2").Style = "Normal"

# --- Rows 3-5: remove the old `ili` example rows, leaving blank styled cells ---
$ws.Range("A3:D5").ClearContents()

# Drop the custom row-height overrides left over from the old (taller) content
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()

# --- Column widths updated to fit the new, shorter content ---
$ws.Columns.Item(1).ColumnWidth = 17.830729166666668
$ws.Columns.Item(2).ColumnWidth = 12.330729166666666
$ws.Columns.Item(3).ColumnWidth = 54.166666666666664
$ws.Columns.Item(4).ColumnWidth = 4.998697916666667

# --- Selection moved ---
[void]$ws.Range("I6").Select()
